$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 5 with 2021 data, following the same pattern as existing rows (2018-2020)
$ws.Range("A5").Value = "2021年"
$ws.Range("B5").Value = 17928
$ws.Range("C5").Value = 5071
$ws.Range("D5").Value = 1182
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 9738
$ws.Range("G5").Value = 238293
$ws.Range("H5").Value = 14869
$ws.Range("I5").Value = 1325
$ws.Range("J5").Value = 6725
$ws.Range("K5").Value = 3460
$ws.Range("L5").Value = 3671
$ws.Range("M5").Value = 245005
$ws.Range("N5").Value = 1112
$ws.Range("O5").Value = 85
$ws.Range("P5").Value = 5241
$ws.Range("Q5").Value = 4558
$ws.Range("R5").Value = 386
$ws.Range("S5").Value = 3575
$ws.Range("T5").Value = 13492
$ws.Range("U5").Value = 832
$ws.Range("V5").Value = 11684
$ws.Range("W5").Value = 106
$ws.Range("X5").Value = 972
$ws.Range("Y5").Value = 602
$ws.Range("Z5").Value = 3981
$ws.Range("AA5").Value = 2547
$ws.Range("AB5").Value = 21526
$ws.Range("AC5").Value = 3653
$ws.Range("AD5").Value = 1054
$ws.Range("AE5").Value = 72
$ws.Range("AF5").Value = 10068
$ws.Range("AG5").Value = 4607
$ws.Range("AH5").Value = 18725
$ws.Range("AI5").Value = 21615
$ws.Range("AJ5").Value = 3446
$ws.Range("AK5").Value = 2909
$ws.Range("AL5").Value = 2731
$ws.Range("AM5").Value = 295
$ws.Range("AN5").Value = 17158
$ws.Range("AO5").Value = 3722
$ws.Range("AP5").Value = 18179
$ws.Range("AQ5").Value = 875
$ws.Range("AR5").Value = 5090
$ws.Range("AS5").Value = 2509
$ws.Range("AT5").Value = 338

# Match the styling of the year label cell (A2:A4) which is bold, centered, bordered
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
